$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

# New user row (row 5): Julián Guardián / Blood / 123456789 / jjuliang.2015@gmail.com / Nivel 3 / Id 4
#
# "123456789" must be stored as text (like the existing "272426799" in C3),
# not auto-converted to a number. Force text via NumberFormat, assign, then
# clear the format again so the cell keeps the default style (no explicit
# style index), matching how the other text-like cells in the sheet look.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "123456789"
$ws.Range("C5").ClearFormats()

$ws.Range("D5").Value = "jjuliang.2015@gmail.com"
$ws.Range("A5").Value = "Julián Guardián"
$ws.Range("B5").Value = "Blood"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 4

$ws.Range("B5").Select()
